$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
# Row 74
$ws1.Range("H74").Value = 4300
$ws1.Range("I74").Value = 4160
$ws1.Range("J74").Value = 5000
$ws1.Range("K74").Value = 4160
$ws1.Range("L74").Value = 5000
$ws1.Range("M74").Value = -3224
$ws1.Range("N74").Value = -6872

# Row 76
$ws1.Range("H76").Value = 3332.3635
$ws1.Range("I76").Value = 2828.5715
$ws1.Range("J76").Value = 5291.5557
$ws1.Range("K76").Value = 2828.5715
$ws1.Range("L76").Value = 5291.5557
$ws1.Range("M76").Value = -2513.5715
$ws1.Range("N76").Value = -5921.5557

# Row 77
$ws1.Range("H77").Value = 4300
$ws1.Range("I77").Value = 4160
$ws1.Range("J77").Value = 5000
$ws1.Range("K77").Value = 20800
$ws1.Range("L77").Value = 25000
$ws1.Range("M77").Value = -16120
$ws1.Range("N77").Value = -34360

# Row 79
$ws1.Range("H79").Value = 3332.3635
$ws1.Range("I79").Value = 2828.5715
$ws1.Range("J79").Value = 5291.5557
$ws1.Range("K79").Value = 2828.5715
$ws1.Range("L79").Value = 5291.5557
$ws1.Range("M79").Value = -1736.5715
$ws1.Range("N79").Value = -7475.5557

# Row 115
$ws1.Range("H115").Value = 641.46155
$ws1.Range("I115").Value = 229.875
$ws1.Range("J115").Value = 1300
$ws1.Range("K115").Value = 689.625
$ws1.Range("L115").Value = 3900
$ws1.Range("M115").Value = 877.375
$ws1.Range("N115").Value = -7034

# Row 116
$ws1.Range("H116").Value = 5591.933
$ws1.Range("I116").Value = 5980
$ws1.Range("J116").Value = 5333.222
$ws1.Range("K116").Value = 5980
$ws1.Range("L116").Value = 5333.222
$ws1.Range("M116").Value = -2538
$ws1.Range("N116").Value = -12217.222

# Row 138
$ws1.Range("H138").Value = 7427.4478
$ws1.Range("I138").Value = 3871.818
$ws1.Range("J138").Value = 8125.875
$ws1.Range("K138").Value = 11615.454
$ws1.Range("L138").Value = 24377.625
$ws1.Range("M138").Value = -6475.454000000002
$ws1.Range("N138").Value = -34657.625

$ws2 = $wb.Worksheets.Item("ARM")
# Row 44
$ws2.Range("H44").Value = 29907.143
$ws2.Range("J44").Value = 29907.143
$ws2.Range("L44").Value = 29907.143
$ws2.Range("N44").Value = -30883.143

# Row 63
$ws2.Range("H63").Value = 10000
$ws2.Range("I63").Value = 0
$ws2.Range("J63").Value = 10000
$ws2.Range("K63").Value = 0
$ws2.Range("L63").Value = 10000
$ws2.Range("M63").ClearContents()
$ws2.Range("N63").Value = -11372

# Row 66
$ws2.Range("H66").Value = 10000
$ws2.Range("I66").Value = 0
$ws2.Range("J66").Value = 10000
$ws2.Range("K66").Value = 0
$ws2.Range("L66").Value = 50000
$ws2.Range("M66").ClearContents()
$ws2.Range("N66").Value = -56864

# Row 102
$ws2.Range("H102").Value = 1546.9231
$ws2.Range("I102").Value = 1283.5454
$ws2.Range("J102").Value = 2995.5
$ws2.Range("K102").Value = 1283.5454
$ws2.Range("L102").Value = 2995.5
$ws2.Range("M102").Value = 338.4546
$ws2.Range("N102").Value = -6239.5

$ws3 = $wb.Worksheets.Item("BSM")
# Row 82
$ws3.Range("H82").Value = 2999.75
$ws3.Range("I82").Value = 2999.75
$ws3.Range("J82").Value = 0
$ws3.Range("K82").Value = 2999.75
$ws3.Range("L82").Value = 0
$ws3.Range("M82").Value = -2616.75
$ws3.Range("N82").ClearContents()

# Row 85
$ws3.Range("H85").Value = 2999.75
$ws3.Range("I85").Value = 2999.75
$ws3.Range("J85").Value = 0
$ws3.Range("K85").Value = 2999.75
$ws3.Range("L85").Value = 0
$ws3.Range("M85").Value = -1673.75
$ws3.Range("N85").ClearContents()

# Row 86
$ws3.Range("H86").Value = 3189.4736
$ws3.Range("I86").Value = 3175
$ws3.Range("J86").Value = 3200
$ws3.Range("K86").Value = 3175
$ws3.Range("L86").Value = 3200
$ws3.Range("M86").Value = -2052
$ws3.Range("N86").Value = -5446

# Row 89
$ws3.Range("H89").Value = 3189.4736
$ws3.Range("I89").Value = 3175
$ws3.Range("J89").Value = 3200
$ws3.Range("K89").Value = 15875
$ws3.Range("L89").Value = 16000
$ws3.Range("M89").Value = -10259
$ws3.Range("N89").Value = -27232

$ws5 = $wb.Worksheets.Item("CUL")
# Row 113
$ws5.Range("H113").Value = 1625.909
$ws5.Range("I113").Value = 2700.2
$ws5.Range("J113").Value = 730.6667
$ws5.Range("K113").Value = 8100.599999999999
$ws5.Range("L113").Value = 2192.0001
$ws5.Range("M113").Value = -5930.599999999999
$ws5.Range("N113").Value = -6532.0001

$ws6 = $wb.Worksheets.Item("GSM")
# Row 70
$ws6.Range("H70").Value = 6314.625
$ws6.Range("I70").Value = 4302.6665
$ws6.Range("J70").Value = 7521.8
$ws6.Range("K70").Value = 4302.6665
$ws6.Range("L70").Value = 7521.8
$ws6.Range("M70").Value = -4032.6665
$ws6.Range("N70").Value = -8061.8

# Row 73
$ws6.Range("H73").Value = 6314.625
$ws6.Range("I73").Value = 4302.6665
$ws6.Range("J73").Value = 7521.8
$ws6.Range("K73").Value = 4302.6665
$ws6.Range("L73").Value = 7521.8
$ws6.Range("M73").Value = -3366.6665
$ws6.Range("N73").Value = -9393.799999999999

# Row 132
$ws6.Range("H132").Value = 2454.0527
$ws6.Range("I132").Value = 1924.2
$ws6.Range("J132").Value = 3042.7778
$ws6.Range("K132").Value = 5772.6
$ws6.Range("L132").Value = 9128.3334
$ws6.Range("M132").Value = -3242.6
$ws6.Range("N132").Value = -14188.3334

$ws7 = $wb.Worksheets.Item("LTW")
# Row 46
$ws7.Range("H46").Value = 640
$ws7.Range("I46").Value = 566.6667
$ws7.Range("J46").Value = 750
$ws7.Range("K46").Value = 566.6667
$ws7.Range("L46").Value = 750
$ws7.Range("M46").Value = -378.6667
$ws7.Range("N46").Value = -1126

# Row 93
$ws7.Range("H93").Value = 587.05554
$ws7.Range("I93").Value = 456.2857
$ws7.Range("J93").Value = 1044.75
$ws7.Range("K93").Value = 456.2857
$ws7.Range("L93").Value = 1044.75
$ws7.Range("M93").Value = 791.7143
$ws7.Range("N93").Value = -3540.75

$ws8 = $wb.Worksheets.Item("WVR")
# Row 81
$ws8.Range("H81").Value = 1633.6666
$ws8.Range("I81").Value = 1000.5
$ws8.Range("J81").Value = 2900
$ws8.Range("K81").Value = 2001
$ws8.Range("L81").Value = 5800
$ws8.Range("M81").Value = -940
$ws8.Range("N81").Value = -7922

# Row 84
$ws8.Range("H84").Value = 1633.6666
$ws8.Range("I84").Value = 1000.5
$ws8.Range("J84").Value = 2900
$ws8.Range("K84").Value = 10005
$ws8.Range("L84").Value = 29000
$ws8.Range("M84").Value = -4701
$ws8.Range("N84").Value = -39608

# Row 132
$ws8.Range("H132").Value = 2279.804
$ws8.Range("I132").Value = 1613.2903
$ws8.Range("J132").Value = 3312.9
$ws8.Range("K132").Value = 4839.8709
$ws8.Range("L132").Value = 9938.700000000001
$ws8.Range("M132").Value = -2309.8709
$ws8.Range("N132").Value = -14998.7
